$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.443826
$ws.Range("H2").Value = 64.331478
$ws.Range("I2").Value = 0.6062978927103765
$ws.Range("J2").Value = 0.6062978927103765
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.379592
$ws.Range("N2").Value = 16.138776
$ws.Range("O2").Value = 0.1347369221116526
$ws.Range("P2").Value = 0.1347369221116527
$ws.Range("Q2").Value = 115.359034798992
$ws.Range("R2").Value = 1038.231313190928
$ws.Range("S2").Value = 0.08169071194657712
$ws.Range("T2").Value = 0.08169071194657716

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.443826
$ws.Range("H3").Value = 64.331478
$ws.Range("I3").Value = 0.6062978927103765
$ws.Range("J3").Value = 0.6062978927103765
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 13.148327
$ws.Range("N3").Value = 39.444981
$ws.Range("O3").Value = 0.3293121691937864
$ws.Range("P3").Value = 0.3293121691937864
$ws.Range("Q3").Value = 281.950436379102
$ws.Range("R3").Value = 2537.553927411918
$ws.Range("S3").Value = 0.1996612742260756
$ws.Range("T3").Value = 0.1996612742260757

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.443826
$ws.Range("H4").Value = 64.331478
$ws.Range("I4").Value = 0.6062978927103765
$ws.Range("J4").Value = 0.6062978927103765
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.33802
$ws.Range("N4").Value = 16.01406
$ws.Range("O4").Value = 0.133695712420281
$ws.Range("P4").Value = 0.133695712420281
$ws.Range("Q4").Value = 114.46757206452
$ws.Range("R4").Value = 1030.20814858068
$ws.Range("S4").Value = 0.08105942870482886
$ws.Range("T4").Value = 0.08105942870482888

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 21.443826
$ws.Range("H5").Value = 64.331478
$ws.Range("I5").Value = 0.6062978927103765
$ws.Range("J5").Value = 0.6062978927103765
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.091788666666666
$ws.Range("N5").Value = 15.275366
$ws.Range("O5").Value = 0.1275286179676195
$ws.Range("P5").Value = 0.1275286179676196
$ws.Range("Q5").Value = 109.187430196772
$ws.Range("R5").Value = 982.686871770948
$ws.Range("S5").Value = 0.07732033233403438
$ws.Range("T5").Value = 0.0773203323340344

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 21.443826
$ws.Range("H6").Value = 64.331478
$ws.Range("I6").Value = 0.6062978927103765
$ws.Range("J6").Value = 0.6062978927103765
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.855806
$ws.Range("N6").Value = 2.567418
$ws.Range("O6").Value = 0.02143446312744256
$ws.Range("P6").Value = 0.02143446312744257
$ws.Range("Q6").Value = 18.351754953756
$ws.Range("R6").Value = 165.165794583804
$ws.Range("S6").Value = 0.01299566982554669
$ws.Range("T6").Value = 0.0129956698255467

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 21.443826
$ws.Range("H7").Value = 64.331478
$ws.Range("I7").Value = 0.6062978927103765
$ws.Range("J7").Value = 0.6062978927103765
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 10.113102
$ws.Range("N7").Value = 30.339306
$ws.Range("O7").Value = 0.2532921151792178
$ws.Range("P7").Value = 0.2532921151792179
$ws.Range("Q7").Value = 216.863599608252
$ws.Range("R7").Value = 1951.772396474268
$ws.Range("S7").Value = 0.1535704756733137
$ws.Range("T7").Value = 0.1535704756733138

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.641794
$ws.Range("H8").Value = 10.925382
$ws.Range("I8").Value = 0.1029672609675761
$ws.Range("J8").Value = 0.1029672609675761
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 5.379592
$ws.Range("N8").Value = 16.138776
$ws.Range("O8").Value = 0.1347369221116526
$ws.Range("P8").Value = 0.1347369221116527
$ws.Range("Q8").Value = 19.591365868048
$ws.Range("R8").Value = 176.322292812432
$ws.Range("S8").Value = 0.01387349182103851
$ws.Range("T8").Value = 0.01387349182103851

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.641794
$ws.Range("H9").Value = 10.925382
$ws.Range("I9").Value = 0.1029672609675761
$ws.Range("J9").Value = 0.1029672609675761
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 13.148327
$ws.Range("N9").Value = 39.444981
$ws.Range("O9").Value = 0.3293121691937864
$ws.Range("P9").Value = 0.3293121691937864
$ws.Range("Q9").Value = 47.88349837863801
$ws.Range("R9").Value = 430.951485407742
$ws.Range("S9").Value = 0.03390837206517516
$ws.Range("T9").Value = 0.03390837206517517

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.641794
$ws.Range("H10").Value = 10.925382
$ws.Range("I10").Value = 0.1029672609675761
$ws.Range("J10").Value = 0.1029672609675761
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.33802
$ws.Range("N10").Value = 16.01406
$ws.Range("O10").Value = 0.133695712420281
$ws.Range("P10").Value = 0.133695712420281
$ws.Range("Q10").Value = 19.43996920788
$ws.Range("R10").Value = 174.95972287092
$ws.Range("S10").Value = 0.01376628131102507
$ws.Range("T10").Value = 0.01376628131102507

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.641794
$ws.Range("H11").Value = 10.925382
$ws.Range("I11").Value = 0.1029672609675761
$ws.Range("J11").Value = 0.1029672609675761
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 5.091788666666666
$ws.Range("N11").Value = 15.275366
$ws.Range("O11").Value = 0.1275286179676195
$ws.Range("P11").Value = 0.1275286179676196
$ws.Range("Q11").Value = 18.54324541553467
$ws.Range("R11").Value = 166.889208739812
$ws.Range("S11").Value = 0.01313127248710619
$ws.Range("T11").Value = 0.0131312724871062

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.641794
$ws.Range("H12").Value = 10.925382
$ws.Range("I12").Value = 0.1029672609675761
$ws.Range("J12").Value = 0.1029672609675761
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.855806
$ws.Range("N12").Value = 2.567418
$ws.Range("O12").Value = 0.02143446312744256
$ws.Range("P12").Value = 0.02143446312744257
$ws.Range("Q12").Value = 3.116669155964
$ws.Range("R12").Value = 28.050022403676
$ws.Range("S12").Value = 0.002207047958543265
$ws.Range("T12").Value = 0.002207047958543266

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.641794
$ws.Range("H13").Value = 10.925382
$ws.Range("I13").Value = 0.1029672609675761
$ws.Range("J13").Value = 0.1029672609675761
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 10.113102
$ws.Range("N13").Value = 30.339306
$ws.Range("O13").Value = 0.2532921151792178
$ws.Range("P13").Value = 0.2532921151792179
$ws.Range("Q13").Value = 36.829834184988
$ws.Range("R13").Value = 331.4685076648921
$ws.Range("S13").Value = 0.02608079532468785
$ws.Range("T13").Value = 0.02608079532468786

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 10.28284533333333
$ws.Range("H14").Value = 30.848536
$ws.Range("I14").Value = 0.2907348463220475
$ws.Range("J14").Value = 0.2907348463220475
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 5.379592
$ws.Range("N14").Value = 16.138776
$ws.Range("O14").Value = 0.1347369221116526
$ws.Range("P14").Value = 0.1347369221116527
$ws.Range("Q14").Value = 55.31751249243732
$ws.Range("R14").Value = 497.857612431936
$ws.Range("S14").Value = 0.03917271834403702
$ws.Range("T14").Value = 0.03917271834403704

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 10.28284533333333
$ws.Range("H15").Value = 30.848536
$ws.Range("I15").Value = 0.2907348463220475
$ws.Range("J15").Value = 0.2907348463220475
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 13.148327
$ws.Range("N15").Value = 39.444981
$ws.Range("O15").Value = 0.3293121691937864
$ws.Range("P15").Value = 0.3293121691937864
$ws.Range("Q15").Value = 135.2022129330907
$ws.Range("R15").Value = 1216.819916397816
$ws.Range("S15").Value = 0.0957425229025356
$ws.Range("T15").Value = 0.09574252290253561

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 10.28284533333333
$ws.Range("H16").Value = 30.848536
$ws.Range("I16").Value = 0.2907348463220475
$ws.Range("J16").Value = 0.2907348463220475
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 5.33802
$ws.Range("N16").Value = 16.01406
$ws.Range("O16").Value = 0.133695712420281
$ws.Range("P16").Value = 0.133695712420281
$ws.Range("Q16").Value = 54.89003404624
$ws.Range("R16").Value = 494.01030641616
$ws.Range("S16").Value = 0.03887000240442705
$ws.Range("T16").Value = 0.03887000240442706

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 10.28284533333333
$ws.Range("H17").Value = 30.848536
$ws.Range("I17").Value = 0.2907348463220475
$ws.Range("J17").Value = 0.2907348463220475
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 5.091788666666666
$ws.Range("N17").Value = 15.275366
$ws.Range("O17").Value = 0.1275286179676195
$ws.Range("P17").Value = 0.1275286179676196
$ws.Range("Q17").Value = 52.35807532935288
$ws.Range("R17").Value = 471.2226779641759
$ws.Range("S17").Value = 0.03707701314647897
$ws.Range("T17").Value = 0.03707701314647899

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 10.28284533333333
$ws.Range("H18").Value = 30.848536
$ws.Range("I18").Value = 0.2907348463220475
$ws.Range("J18").Value = 0.2907348463220475
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 0.855806
$ws.Range("N18").Value = 2.567418
$ws.Range("O18").Value = 0.02143446312744256
$ws.Range("P18").Value = 0.02143446312744257
$ws.Range("Q18").Value = 8.800120733338666
$ws.Range("R18").Value = 79.201086600048
$ws.Range("S18").Value = 0.006231745343352609
$ws.Range("T18").Value = 0.006231745343352611

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 10.28284533333333
$ws.Range("H19").Value = 30.848536
$ws.Range("I19").Value = 0.2907348463220475
$ws.Range("J19").Value = 0.2907348463220475
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 10.113102
$ws.Range("N19").Value = 30.339306
$ws.Range("O19").Value = 0.2532921151792178
$ws.Range("P19").Value = 0.2532921151792179
$ws.Range("Q19").Value = 103.991463706224
$ws.Range("R19").Value = 935.923173356016
$ws.Range("S19").Value = 0.07364084418121625
$ws.Range("T19").Value = 0.07364084418121629
